$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Room Registration")

# Rename the sheet; Excel automatically updates the _FilterDatabase
# defined name's sheet-name qualifier to match.
$ws.Name = "Room Booking"

# Move the active selection on the sheet from G8 to B1.
$ws.Range("B1").Select()
